# Edit the Field Illumination schema sheets:
#  - FieldIlluminationInput: remove "center_threshold" column, shift remaining columns left
#  - FieldIlluminationOutput: rename "profile_rois"->"roi_profiles", "corner_rois"->"roi_corners",
#                              "center_of_illumination"->"roi_centroids_weighted"
#  - FieldIlluminationKeyValues: split "center_fraction" into "center_region_intensity_fraction" and
#                                 "center_region_area_fraction"; add relative-position / relative-distance
#                                 columns for centroid_weighted, centroid, add centroid_fitted columns,
#                                 reorder max_intensity_pos_x/y (y first) and add their relative/distance columns

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# FieldIlluminationInput (was A1:G1, now A1:F1 - center_threshold removed)
# ---------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("FieldIlluminationInput")

$inputHeaders = @(
    "field_illumination_image",
    "bit_depth",
    "saturation_threshold",
    "corner_fraction",
    "sigma",
    "intensity_map_size"
)

for ($i = 0; $i -lt $inputHeaders.Length; $i++) {
    $wsInput.Cells.Item(1, $i + 1).Value = $inputHeaders[$i]
}
$wsInput.Cells.Item(1, 7).ClearContents()

# ---------------------------------------------------------------
# FieldIlluminationOutput (still A1:F1, 3 of the cells renamed)
# ---------------------------------------------------------------
$wsOutput = $wb.Worksheets.Item("FieldIlluminationOutput")

$outputHeaders = @(
    "key_values",
    "intensity_profiles",
    "intensity_map",
    "roi_profiles",
    "roi_corners",
    "roi_centroids_weighted"
)

for ($i = 0; $i -lt $outputHeaders.Length; $i++) {
    $wsOutput.Cells.Item(1, $i + 1).Value = $outputHeaders[$i]
}

# ---------------------------------------------------------------
# FieldIlluminationKeyValues (was A1:AK1, now A1:AZ1)
# ---------------------------------------------------------------
$wsKeyValues = $wb.Worksheets.Item("FieldIlluminationKeyValues")

$keyValueHeaders = @(
    "channel",
    "center_region_intensity_fraction",
    "center_region_area_fraction",
    "centroid_weighted_y",
    "centroid_weighted_y_relative",
    "centroid_weighted_x",
    "centroid_weighted_x_relative",
    "centroid_weighted_distance_relative",
    "centroid_y",
    "centroid_y_relative",
    "centroid_x",
    "centroid_x_relative",
    "centroid_distance_relative",
    "centroid_fitted_y",
    "centroid_fitted_y_relative",
    "centroid_fitted_x",
    "centroid_fitted_x_relative",
    "centroid_fitted_distance_relative",
    "max_intensity",
    "max_intensity_pos_y",
    "max_intensity_pos_y_relative",
    "max_intensity_pos_x",
    "max_intensity_pos_x_relative",
    "max_intensity_distance_relative",
    "top_left_intensity_mean",
    "top_left_intensity_ratio",
    "top_center_intensity_mean",
    "top_center_intensity_ratio",
    "top_right_intensity_mean",
    "top_right_intensity_ratio",
    "middle_left_intensity_mean",
    "middle_left_intensity_ratio",
    "middle_center_intensity_mean",
    "middle_center_intensity_ratio",
    "middle_right_intensity_mean",
    "middle_right_intensity_ratio",
    "bottom_left_intensity_mean",
    "bottom_left_intensity_ratio",
    "bottom_center_intensity_mean",
    "bottom_center_intensity_ratio",
    "bottom_right_intensity_mean",
    "bottom_right_intensity_ratio",
    "decile_0",
    "decile_1",
    "decile_2",
    "decile_3",
    "decile_4",
    "decile_5",
    "decile_6",
    "decile_7",
    "decile_8",
    "decile_9"
)

for ($i = 0; $i -lt $keyValueHeaders.Length; $i++) {
    $wsKeyValues.Cells.Item(1, $i + 1).Value = $keyValueHeaders[$i]
}
